# "Agendas for meetings 24 and 25"
# Mark the agenda checklist: Ankita (column C) submitted her agenda for
# "Easter Week 5" (row 11) and "Summer Week 1" (row 12) -- both cells were
# previously empty (red "missing" fill) and should become a checkmark,
# matching the formatting already used for the rest of column C / the
# neighboring "done" cells (e.g. D11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$checkMark = [string][char]0x2714

# Copy the existing checkmark cell's formatting (font "Zapf Dingbats",
# green fill, etc.) onto C11 and C12, then write the checkmark value.
$ws.Range("D11").Copy()
$ws.Range("C11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C12").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C11").Value = $checkMark
$ws.Range("C12").Value = $checkMark

# Reflect the author's final selection on the sheet.
$ws.Range("C12").Select()
